$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "<tab>Date :03-09-23" run from the first paragraph, leaving
#    only "Practical-4".
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$tabPos = $p1.Range.Text.IndexOf([char]9)
if ($tabPos -ge 0) {
    $delStart = $p1.Range.Start + $tabPos
    $delEnd = $p1.Range.End - 1
    $delRng = $d.Range($delStart, $delEnd)
    $delRng.Delete()
}

# ---------------------------------------------------------------------------
# 2) Global font change: "Agency FB" -> "Arial Narrow" everywhere (runs and
#    paragraph marks).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    $rng.Font.Name = "Arial Narrow"
    $after = $rng.Font.Name
    if ($after -ne "Arial Narrow") {
        # Paragraph has no run content at all (pure paragraph-mark-only
        # paragraph) - Range.Font.Name does not reach the pPr/rPr in that
        # case. Insert a scratch character so the paragraph temporarily
        # has a real run, (re)apply the font to the whole paragraph range
        # (which now also touches the paragraph mark), then delete the
        # scratch character again.
        $rng.InsertBefore("x")
        $rng2 = $p.Range
        $rng2.Font.Name = "Arial Narrow"
        $scratchStart = $p.Range.Start
        $scratchRng = $d.Range($scratchStart, $scratchStart + 1)
        $scratchRng.Delete()
    }
}

# ---------------------------------------------------------------------------
# 3) Paragraph "Before we can do anything ..." -> italic, size 15pt (sz/szCs
#    30 half-points), font already handled above.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Before we can do anything*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4FD98B13" w14:textId="77777777" w:rsidR="0062077C" w:rsidRPr="007678D8" w:rsidRDefault="0062077C" w:rsidP="0062077C"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Times New Roman" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="30"/><w:szCs w:val="30"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="007678D8"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Times New Roman" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:i/><w:iCs/><w:color w:val="212529"/><w:kern w:val="0"/><w:sz w:val="30"/><w:szCs w:val="30"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Before we can do anything such as finding out the contents of any files or folders, we need to know what exists in the first place. This can be done using the "ls" command (short for listing)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.Range.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 4) The picture-only paragraph right after it gains explicit sz/szCs = 24
#    (12pt), on both the paragraph mark and the run that hosts the drawing.
# ---------------------------------------------------------------------------
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Before we can do anything*") {
        $target2 = $d.Paragraphs.Item($i + 1)
        break
    }
}
if ($target2 -ne $null) {
    $target2.Range.Font.Size = 12
    $target2.Range.Font.SizeBi = 12
}

Write-Output "edit-complete"
